# Updated data sheet NGC-1345
#
# Adds a new "Status In B" column (column I) with "False" values to the
# TC_55046 and TC_55052 sheets, and moves the active-tab / selection
# around the three sheets so that TC_55056 ends up the active tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # TC_55046
$ws2 = $wb.Worksheets.Item(2)   # TC_55052
$ws3 = $wb.Worksheets.Item(3)   # TC_55056

# ---------------------------------------------------------------------
# TC_55046 (sheet1): new column I, rows 7-11
# ---------------------------------------------------------------------

# I7 = header "Status In B", formatted like the other header cells (H7)
$ws1.Range("I7").Formula = "'Status In B"
$ws1.Range("H7").Copy()
$ws1.Range("I7").PasteSpecial(-4122)   # xlPasteFormats

# I8:I11 = "False", formatted like the other value cells (H8)
$ws1.Range("I8").Formula = "'False"
$ws1.Range("I9").Formula = "'False"
$ws1.Range("I10").Formula = "'False"
$ws1.Range("I11").Formula = "'False"
$ws1.Range("H8").Copy()
$ws1.Range("I8:I11").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# TC_55052 (sheet2): existing column I (rows 7-9) gets real content
# ---------------------------------------------------------------------

$ws2.Range("I7").Formula = "'Status In B"
$ws2.Range("H7").Copy()
$ws2.Range("I7").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("I8").Formula = "'False"
$ws2.Range("I9").Formula = "'False"
$ws2.Range("A8").Copy()
$ws2.Range("I8:I9").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Selections / active tab
#
# Originally TC_55052 (index 1) is tabSelected; the edit moves that to
# TC_55056 (index 2). Touch sheets in the order that leaves TC_55056
# activated last.
# ---------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("I7:I9").Select()

$ws2.Activate()
$ws2.Range("I7:I9").Select()

$ws3.Activate()
$ws3.Range("H7").Select()
